$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '64.054.43'
$ws.Cells.Item(2, 5).Value = '  +1.99%  '

$ws.Cells.Item(3, 4).Value = '3.442.97'
$ws.Cells.Item(3, 5).Value = '  +3.16%  '

$ws.Cells.Item(4, 5).Value = '  +0.02%  '

$ws.Cells.Item(5, 4).Value = "'574.86"
$ws.Cells.Item(5, 5).Value = '  +3.59%  '

$ws.Cells.Item(6, 4).Value = "'157.07"
$ws.Cells.Item(6, 5).Value = '  +3.32%  '

$ws.Cells.Item(7, 5).Value = '  +0.09%  '

$ws.Cells.Item(8, 4).Value = '3.444.47'
$ws.Cells.Item(8, 5).Value = '  +3.00%  '

$ws.Cells.Item(9, 4).Value = "'0.549"
$ws.Cells.Item(9, 5).Value = '  +3.61%  '

$ws.Cells.Item(10, 4).Value = "'7.52"
$ws.Cells.Item(10, 5).Value = '  +0.13%  '

$ws.Cells.Item(11, 5).Value = '  +4.84%  '

$ws.Cells.Item(12, 5).Value = '  +1.29%  '

$ws.Cells.Item(13, 4).Value = '4.039.43'
$ws.Cells.Item(13, 5).Value = '  +3.41%  '

$ws.Cells.Item(14, 5).Value = '  -2.17%  '

$ws.Cells.Item(15, 5).Value = '  +7.60%  '

$ws.Cells.Item(16, 4).Value = "'27.41"
$ws.Cells.Item(16, 5).Value = '  +2.28%  '

$ws.Cells.Item(17, 4).Value = '64.111.14'
$ws.Cells.Item(17, 5).Value = '  +2.11%  '

$ws.Cells.Item(18, 4).Value = '3.461.78'
$ws.Cells.Item(18, 5).Value = '  +3.34%  '

$ws.Cells.Item(19, 5).Value = '  -1.37%  '

$ws.Cells.Item(20, 4).Value = "'14.33"
$ws.Cells.Item(20, 5).Value = '  +4.39%  '

$ws.Cells.Item(21, 4).Value = "'392.11"
$ws.Cells.Item(21, 5).Value = '  +0.72%  '

$ws.Cells.Item(22, 4).Value = "'8.39"
$ws.Cells.Item(22, 5).Value = '  -0.69%  '

$ws.Cells.Item(23, 2).Value = 'Polygon'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(23, 4).Value = "'0.543"
$ws.Cells.Item(23, 5).Value = '  +0.87%  '

$ws.Cells.Item(24, 2).Value = 'Litecoin'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(24, 4).Value = "'72.89"
$ws.Cells.Item(24, 5).Value = '  +3.20%  '

$ws.Cells.Item(25, 4).Value = "'0.995"
$ws.Cells.Item(25, 5).Value = '  -0.48%  '

$ws.Cells.Item(26, 5).Value = '  +26.20%  '

$ws.Cells.Item(27, 4).Value = "'9.67"
$ws.Cells.Item(27, 5).Value = '  +9.81%  '

$ws.Cells.Item(28, 4).Value = "'0.178"
$ws.Cells.Item(28, 5).Value = '  +0.01%  '

$ws.Cells.Item(29, 4).Value = "'0.999"
$ws.Cells.Item(29, 5).Value = '  -0.21%  '

$ws.Cells.Item(30, 4).Value = "'6.09"
$ws.Cells.Item(30, 5).Value = '  +9.23%  '

$ws.Cells.Item(31, 5).Value = '  +2.73%  '

$ws.Cells.Item(32, 5).Value = '  +6.26%  '

$ws.Cells.Item(33, 4).Value = "'23.65"
$ws.Cells.Item(33, 5).Value = '  +2.86%  '

$ws.Cells.Item(34, 4).Value = "'6.53"
$ws.Cells.Item(34, 5).Value = '  +0.85%  '

$ws.Cells.Item(35, 5).Value = '  +0.02%  '

$ws.Cells.Item(36, 4).Value = "'7.04"
$ws.Cells.Item(36, 5).Value = '  +5.54%  '

$ws.Cells.Item(37, 4).Value = "'160.69"
$ws.Cells.Item(37, 5).Value = '  -0.11%  '

$ws.Cells.Item(38, 4).Value = "'1.47"
$ws.Cells.Item(38, 5).Value = '  -0.91%  '

$ws.Cells.Item(39, 4).Value = "'0.0782"
$ws.Cells.Item(39, 5).Value = '  +6.45%  '

$ws.Cells.Item(40, 2).Value = 'Stacks'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(40, 4).Value = "'1.85"
$ws.Cells.Item(40, 5).Value = '  -1.48%  '

$ws.Cells.Item(41, 4).Value = '2.929.64'
$ws.Cells.Item(41, 5).Value = '  +2.86%  '

$ws.Cells.Item(42, 2).Value = 'EnergySwap'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(42, 4).Value = "'27.23"
$ws.Cells.Item(42, 5).Value = '  -0.16%  '

$ws.Cells.Item(43, 5).Value = '  +2.37%  '

$ws.Cells.Item(44, 5).Value = '  +2.68%  '

$ws.Cells.Item(45, 4).Value = "'0.772"
$ws.Cells.Item(45, 5).Value = '  +3.11%  '

$ws.Cells.Item(46, 4).Value = "'41.90"
$ws.Cells.Item(46, 5).Value = '  +2.82%  '

$ws.Cells.Item(47, 4).Value = "'23.73"
$ws.Cells.Item(47, 5).Value = '  +7.92%  '

$ws.Cells.Item(48, 5).Value = '  +4.43%  '

$ws.Cells.Item(49, 5).Value = '  +23.69%  '

$ws.Cells.Item(50, 5).Value = '  +6.92%  '

$ws.Cells.Item(51, 4).Value = "'6.54"
$ws.Cells.Item(51, 5).Value = '  +4.22%  '
